$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Feature 2 file"
$ws.Range("A2").Select()
